$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 484-485, shifting the existing rows 484-588 down to 486-590.
$ws.Rows("484:485").Insert()

# --- New row 484 ---
$ws.Range("A484").Value = 7
$ws.Range("B484").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C484").Value = "Ñuble"
$ws.Range("D484").Value = 45204
$ws.Range("E484").Value = 16
$ws.Range("F484").Value = 100112009
$ws.Range("G484").Value = "Acelga"
$ws.Range("H484").Value = "Sin especificar"
$ws.Range("I484").Value = "Primera"
$ws.Range("J484").Value = 300
$ws.Range("K484").Value = 700
$ws.Range("L484").Value = 700
$ws.Range("M484").Value = 700
$ws.Range("N484").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O484").Value = "Provincia de Diguillín"
$ws.Range("P484").Value = 700
$ws.Range("Q484").Value = 1
$ws.Range("R484").Value = "Hortaliza"

# --- New row 485 ---
$ws.Range("A485").Value = 7
$ws.Range("B485").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C485").Value = "Ñuble"
$ws.Range("D485").Value = 45204
$ws.Range("E485").Value = 16
$ws.Range("F485").Value = 100112009
$ws.Range("G485").Value = "Acelga"
$ws.Range("H485").Value = "Sin especificar"
$ws.Range("I485").Value = "Segunda"
$ws.Range("J485").Value = 300
$ws.Range("K485").Value = 500
$ws.Range("L485").Value = 500
$ws.Range("M485").Value = 500
$ws.Range("N485").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O485").Value = "Provincia de Diguillín"
$ws.Range("P485").Value = 500
$ws.Range("Q485").Value = 1
$ws.Range("R485").Value = "Hortaliza"
